$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2: iaest-dimension:tipo-de-hogar -> iaest-measure:tipo-de-hogar
$ws.Range("B2").Value = "iaest-measure:tipo-de-hogar"

# Update B3: dim -> medida
$ws.Range("B3").Value = "medida"

# Update B4: skos:Concept -> xsd:int
$ws.Range("B4").Value = "xsd:int"

# Delete row 5 (contains mapping-tipo-de-hogar.xlsx in B5)
$ws.Rows("5").Delete()
